$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) homePage.searchBarAriaLabel row (row 5): the en-US translation was
#    missing. C5 currently carries the "empty" cell style (no number format);
#    match the filled-cell style used by its siblings (e.g. B5) before
#    writing the value, so it looks identical to every other populated cell.
# ---------------------------------------------------------------------------
$ws.Range("B5").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = "Type to Search"

# ---------------------------------------------------------------------------
# 2) Insert a new translation row (feedback.subtitle) right after the
#    feedback.hero row (row 29), i.e. it becomes the new row 30. Rather than
#    using Rows.Insert (which in this host synthesises brand-new,
#    unreferenced cell styles), shift the existing rows 30..46 down into
#    31..47 by copying each row (values + formats + row height) into the row
#    below it, working bottom-to-top so nothing gets clobbered before it's
#    copied. Row 30 is then overwritten with the new feedback.subtitle
#    content (row 29, feedback.hero, is left untouched).
# ---------------------------------------------------------------------------
for ($r = 46; $r -ge 30; $r--) {
    $src = $r
    $dst = $r + 1
    $ws.Range("A$src`:E$src").Copy()
    $ws.Range("A$dst`:E$dst").PasteSpecial(-4122)
    $ws.Range("A$src`:E$src").Copy()
    $ws.Range("A$dst`:E$dst").PasteSpecial(-4163)
    $ws.Rows.Item($dst).RowHeight = $ws.Rows.Item($src).RowHeight
}

$ws.Range("A30").Value = "feedback.subtitle"
$ws.Range("B30").Value = "你可以畅所欲言"
$ws.Range("C30").Value = "New app request, bug report, or anything you want to tell us."

# ---------------------------------------------------------------------------
# 3) Append three new translation rows at the bottom of the table
#    (app.decision.addOption / savePreset / currentOption), reusing the
#    format of the last existing row (47, originally aboutPage.meta.title).
# ---------------------------------------------------------------------------
$ws.Range("A47:E47").Copy()
$ws.Range("A48:E48").PasteSpecial(-4122)
$ws.Range("A48").Value = "app.decision.addOption"
$ws.Range("B48").Value = "添加选项"
$ws.Range("C48").Value = "Add Option"
$ws.Rows.Item(48).RowHeight = 20.1

$ws.Range("A48:E48").Copy()
$ws.Range("A49:E49").PasteSpecial(-4122)
$ws.Range("A49").Value = "app.decision.savePreset"
$ws.Range("B49").Value = "保存预设"
$ws.Range("C49").Value = "Save Preset"
$ws.Rows.Item(49).RowHeight = 20.1

$ws.Range("A49:E49").Copy()
$ws.Range("A50:E50").PasteSpecial(-4122)
$ws.Range("A50").Value = "app.decision.currentOption"
$ws.Range("B50").Value = "当前备选项"
$ws.Range("C50").Value = "Current Options"
$ws.Rows.Item(50).RowHeight = 20.1
